$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the "Timed Out" (TO) marking for QuantLibAddin (row 6) into the
# VC12 (2013) column group (R:Y), matching the existing VC11 (2012) group (J:Q).
$ws.Range("J6").Copy() | Out-Null
$ws.Range("R6:Y6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6:Y6").Value = "TO"

# Update the notes text to mention both VC11 and VC12 now that VC12 also
# exhibits the same batch-build hang.
$ws.Range("A20").Value = "VC11 / VC12 - batch build / rebuild of QuantLibAddin - kick it off in the evening - following morning it's still running."
